$wb = $excel.ActiveWorkbook

# --- Sheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Name = "summ37847484"
$ws.Range("B2").Value = [double]"8471.484023576697"
$ws.Range("C2").Value = [double]"0.007789454624852856"
$ws.Range("B3").Value = [double]"1202.91681436473"
$ws.Range("C3").Value = [double]"0.3237179624293176"
$ws.Range("B4").Value = [double]"-159.1548727820026"
$ws.Range("C4").Value = [double]"0.9039989314753765"
$ws.Range("B5").Value = [double]"1343.50409074776"
$ws.Range("C5").Value = [double]"0.001904795805325258"
$ws.Range("B6").Value = [double]"-423.4640039928167"
$ws.Range("C6").Value = [double]"0.340237286704573"
$ws.Range("B7").Value = [double]"428.4818319971113"
$ws.Range("C7").Value = [double]"0.4024079324924151"
$ws.Range("B8").Value = [double]"105.4050962795262"
$ws.Range("C8").Value = [double]"0.8153166625886411"
$ws.Range("B9").Value = [double]"-51.99473862204786"
$ws.Range("C9").Value = [double]"0.7270516013796073"
$ws.Range("B10").Value = [double]"-1792.415801216938"
$ws.Range("C10").Value = [double]"1.027174142114593e-07"
$ws.Range("B11").Value = [double]"-68.6323411082019"
$ws.Range("C11").Value = [double]"1.069133753815934e-05"
$ws.Range("B12").Value = [double]"332.8034421178901"
$ws.Range("C12").Value = [double]"0.02741358790951412"
$ws.Range("B13").Value = [double]"518.7276375599293"
$ws.Range("C13").Value = [double]"8.499176814484142e-09"
$ws.Range("B14").Value = [double]"0.04390360215938952"
$ws.Range("C14").Value = [double]"0.6018682252884338"
$ws.Range("B15").Value = [double]"0.0001544284310010299"
$ws.Range("C15").Value = [double]"0.1134588463734251"
$ws.Range("B16").Value = [double]"-25.57114255722337"
$ws.Range("C16").Value = [double]"0.2212344545137763"
$ws.Range("B17").Value = [double]"8.238285150004652"
$ws.Range("C17").Value = [double]"0.7035076161521496"
$ws.Range("B18").Value = [double]"-1121.221914560181"
$ws.Range("C18").Value = [double]"0.5063703655031857"
$ws.Range("B19").Value = [double]"175.2752195898204"
$ws.Range("C19").Value = [double]"0.9317445835755375"
$ws.Rows.Item(20).Delete()

# --- Sheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Name = "summ38128476"
$ws.Range("B2").Value = [double]"9388.195090780291"
$ws.Range("C2").Value = [double]"0.002701674475805872"
$ws.Range("B3").Value = [double]"-451.400434536899"
$ws.Range("C3").Value = [double]"0.7067964718492141"
$ws.Range("B4").Value = [double]"-888.0985406152272"
$ws.Range("C4").Value = [double]"0.4807566086778076"
$ws.Range("B5").Value = [double]"1299.012275283976"
$ws.Range("C5").Value = [double]"0.002859400648191788"
$ws.Range("B6").Value = [double]"-126.0076330132374"
$ws.Range("C6").Value = [double]"0.777457866470011"
$ws.Range("B7").Value = [double]"556.2778009023741"
$ws.Range("C7").Value = [double]"0.2795875388886331"
$ws.Range("B8").Value = [double]"-119.2040214654118"
$ws.Range("C8").Value = [double]"0.7915469735947752"
$ws.Range("B9").Value = [double]"-181.3939428661145"
$ws.Range("C9").Value = [double]"0.221080816497548"
$ws.Range("B10").Value = [double]"-1565.124921901055"
$ws.Range("C10").Value = [double]"3.651503797419026e-06"
$ws.Range("B11").Value = [double]"-67.79046610848874"
$ws.Range("C11").Value = [double]"1.578533431066389e-05"
$ws.Range("B12").Value = [double]"371.4342648783377"
$ws.Range("C12").Value = [double]"0.01491517339252044"
$ws.Range("B13").Value = [double]"536.1449361066989"
$ws.Range("C13").Value = [double]"3.195449424177433e-09"
$ws.Range("B14").Value = [double]"0.04536509919653772"
$ws.Range("C14").Value = [double]"0.5947825596164378"
$ws.Range("B15").Value = [double]"0.0001478251306894226"
$ws.Range("C15").Value = [double]"0.1336174544394857"
$ws.Range("B16").Value = [double]"-23.13032235971275"
$ws.Range("C16").Value = [double]"0.2697829529106897"
$ws.Range("B17").Value = [double]"0.1829825735360178"
$ws.Range("C17").Value = [double]"0.9931736160825645"
$ws.Range("B18").Value = [double]"-1139.687690707527"
$ws.Range("C18").Value = [double]"0.4966918806460892"
$ws.Range("B19").Value = [double]"-162.0825847955971"
$ws.Range("C19").Value = [double]"0.9371716695104371"
$ws.Rows.Item(20).Delete()

# --- Sheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Name = "summ38390602"
$ws.Range("B2").Value = [double]"8014.076755941342"
$ws.Range("C2").Value = [double]"0.01116143965712045"
$ws.Range("B3").Value = [double]"-563.8481305871351"
$ws.Range("C3").Value = [double]"0.6174072916057645"
$ws.Range("B4").Value = [double]"-650.8990563248137"
$ws.Range("C4").Value = [double]"0.6304849457778881"
$ws.Range("B5").Value = [double]"1043.496637496612"
$ws.Range("C5").Value = [double]"0.01647454936231699"
$ws.Range("B6").Value = [double]"79.77899430493127"
$ws.Range("C6").Value = [double]"0.8571927588910685"
$ws.Range("B7").Value = [double]"677.344764391248"
$ws.Range("C7").Value = [double]"0.1854443513866291"
$ws.Range("B8").Value = [double]"132.3491331418018"
$ws.Range("C8").Value = [double]"0.7680359839300567"
$ws.Range("B9").Value = [double]"-106.6847499392817"
$ws.Range("C9").Value = [double]"0.4687200899913418"
$ws.Range("B10").Value = [double]"-1612.752160511142"
$ws.Range("C10").Value = [double]"1.669533465519624e-06"
$ws.Range("B11").Value = [double]"-67.32304504868026"
$ws.Range("C11").Value = [double]"1.789329378116316e-05"
$ws.Range("B12").Value = [double]"459.6000328632737"
$ws.Range("C12").Value = [double]"0.002393627830145819"
$ws.Range("B13").Value = [double]"552.3509917669783"
$ws.Range("C13").Value = [double]"9.099364163634893e-10"
$ws.Range("B14").Value = [double]"0.04531736300271553"
$ws.Range("C14").Value = [double]"0.5945937378289288"
$ws.Range("B15").Value = [double]"0.0001183842392713263"
$ws.Range("C15").Value = [double]"0.2271738556548191"
$ws.Range("B16").Value = [double]"-9.540648950135813"
$ws.Range("C16").Value = [double]"0.6447417735473917"
$ws.Range("B17").Value = [double]"2.601363622802335"
$ws.Range("C17").Value = [double]"0.9039219068053328"
$ws.Range("B18").Value = [double]"-129.2872323760942"
$ws.Range("C18").Value = [double]"0.9390024977104197"
$ws.Range("B19").Value = [double]"-316.866903504152"
$ws.Range("C19").Value = [double]"0.8762734446991882"
$ws.Rows.Item(20).Delete()

# --- Sheet 4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Name = "summ38651166"
$ws.Range("B2").Value = [double]"8192.39428238071"
$ws.Range("C2").Value = [double]"0.009495924628474582"
$ws.Range("B3").Value = [double]"-72.37535507570141"
$ws.Range("C3").Value = [double]"0.9501452931032324"
$ws.Range("B4").Value = [double]"-820.0089836496758"
$ws.Range("C4").Value = [double]"0.530825412281638"
$ws.Range("B5").Value = [double]"1118.668414312749"
$ws.Range("C5").Value = [double]"0.009376123052208415"
$ws.Range("B6").Value = [double]"-216.6081757430255"
$ws.Range("C6").Value = [double]"0.6212429240950313"
$ws.Range("B7").Value = [double]"769.6073206336387"
$ws.Range("C7").Value = [double]"0.131556930117415"
$ws.Range("B8").Value = [double]"39.37460614092288"
$ws.Range("C8").Value = [double]"0.9297525899792032"
$ws.Range("B9").Value = [double]"-104.6030547941709"
$ws.Range("C9").Value = [double]"0.4751171598847709"
$ws.Range("B10").Value = [double]"-1587.191947043571"
$ws.Range("C10").Value = [double]"1.926146621029515e-06"
$ws.Range("B11").Value = [double]"-63.18524486528561"
$ws.Range("C11").Value = [double]"5.706743621100948e-05"
$ws.Range("B12").Value = [double]"347.0568621603508"
$ws.Range("C12").Value = [double]"0.02074430133079653"
$ws.Range("B13").Value = [double]"576.2717069552249"
$ws.Range("C13").Value = [double]"1.46983807668839e-10"
$ws.Range("B14").Value = [double]"0.121080034424218"
$ws.Range("C14").Value = [double]"0.1605366277659063"
$ws.Range("B15").Value = [double]"8.004689215819033e-05"
$ws.Range("C15").Value = [double]"0.4242954090772512"
$ws.Range("B16").Value = [double]"-27.48722775854593"
$ws.Range("C16").Value = [double]"0.1832986539178213"
$ws.Range("B17").Value = [double]"4.703640735839851"
$ws.Range("C17").Value = [double]"0.8276361923847593"
$ws.Range("B18").Value = [double]"116.589060714145"
$ws.Range("C18").Value = [double]"0.9454799807301836"
$ws.Range("B19").Value = [double]"-143.1389313928144"
$ws.Range("C19").Value = [double]"0.9439395806468924"
$ws.Rows.Item(20).Delete()

# --- Sheet 5 ---
$ws = $wb.Worksheets.Item(5)
$ws.Name = "summ38909914"
$ws.Range("B2").Value = [double]"8751.169819745388"
$ws.Range("C2").Value = [double]"0.005707968929948037"
$ws.Range("B3").Value = [double]"187.0190964763813"
$ws.Range("C3").Value = [double]"0.8697382903927412"
$ws.Range("B4").Value = [double]"-383.0767254129623"
$ws.Range("C4").Value = [double]"0.7723595915410544"
$ws.Range("B5").Value = [double]"1302.819393166094"
$ws.Range("C5").Value = [double]"0.002576685220418289"
$ws.Range("B6").Value = [double]"-99.72726620728145"
$ws.Range("C6").Value = [double]"0.8225381865168346"
$ws.Range("B7").Value = [double]"559.2214746220933"
$ws.Range("C7").Value = [double]"0.2733203423510155"
$ws.Range("B8").Value = [double]"331.6639648847795"
$ws.Range("C8").Value = [double]"0.4613817336049707"
$ws.Range("B9").Value = [double]"-91.18560370361313"
$ws.Range("C9").Value = [double]"0.5365185263878008"
$ws.Range("B10").Value = [double]"-1620.959411719839"
$ws.Range("C10").Value = [double]"1.485061836070658e-06"
$ws.Range("B11").Value = [double]"-59.63018674954142"
$ws.Range("C11").Value = [double]"0.0001410393753195"
$ws.Range("B12").Value = [double]"502.1841368984371"
$ws.Range("C12").Value = [double]"0.0009302229800476914"
$ws.Range("B13").Value = [double]"499.7322944245885"
$ws.Range("C13").Value = [double]"2.22531795431986e-08"
$ws.Range("B14").Value = [double]"0.02742375724374134"
$ws.Range("C14").Value = [double]"0.7489304608516021"
$ws.Range("B15").Value = [double]"0.0002218889529390399"
$ws.Range("C15").Value = [double]"0.02390425446256137"
$ws.Range("B16").Value = [double]"-26.33034151363258"
$ws.Range("C16").Value = [double]"0.2057273720078962"
$ws.Range("B17").Value = [double]"-0.3716400372932522"
$ws.Range("C17").Value = [double]"0.9862620756290588"
$ws.Range("B18").Value = [double]"-1103.386608357473"
$ws.Range("C18").Value = [double]"0.5112249857023348"
$ws.Range("B19").Value = [double]"-993.1501381763181"
$ws.Range("C19").Value = [double]"0.6279637282058537"
$ws.Rows.Item(20).Delete()

# --- Sheet 6 ---
$ws = $wb.Worksheets.Item(6)
$ws.Name = "summ39173648"
$ws.Range("B2").Value = [double]"9708.8070372699"
$ws.Range("C2").Value = [double]"0.002257743551191928"
$ws.Range("B3").Value = [double]"-115.7910962592538"
$ws.Range("C3").Value = [double]"0.920109314603994"
$ws.Range("B4").Value = [double]"-500.0978809379999"
$ws.Range("C4").Value = [double]"0.7107119889942686"
$ws.Range("B5").Value = [double]"1128.381252950388"
$ws.Range("C5").Value = [double]"0.01021292020873957"
$ws.Range("B6").Value = [double]"-248.8903530957777"
$ws.Range("C6").Value = [double]"0.5779417563093279"
$ws.Range("B7").Value = [double]"692.9108110197456"
$ws.Range("C7").Value = [double]"0.1802948517163949"
$ws.Range("B8").Value = [double]"-246.9220310302006"
$ws.Range("C8").Value = [double]"0.5853953893056214"
$ws.Range("B9").Value = [double]"-88.64083544199711"
$ws.Range("C9").Value = [double]"0.5505927265654222"
$ws.Range("B10").Value = [double]"-1644.589100190785"
$ws.Range("C10").Value = [double]"1.216321880118561e-06"
$ws.Range("B11").Value = [double]"-68.3100160627849"
$ws.Range("C11").Value = [double]"1.649784719354718e-05"
$ws.Range("B12").Value = [double]"355.6436995071334"
$ws.Range("C12").Value = [double]"0.02043937368963048"
$ws.Range("B13").Value = [double]"522.8118850353435"
$ws.Range("C13").Value = [double]"6.736089196696919e-09"
$ws.Range("B14").Value = [double]"0.007280426873631491"
$ws.Range("C14").Value = [double]"0.9332724953650509"
$ws.Range("B15").Value = [double]"0.0002064539311885584"
$ws.Range("C15").Value = [double]"0.04448440589477411"
$ws.Range("B16").Value = [double]"-16.55787813571613"
$ws.Range("C16").Value = [double]"0.4316486823126681"
$ws.Range("B17").Value = [double]"-0.9311830723898886"
$ws.Range("C17").Value = [double]"0.9657668147682059"
$ws.Range("B18").Value = [double]"-1809.279104195648"
$ws.Range("C18").Value = [double]"0.2916981662285077"
$ws.Range("B19").Value = [double]"-791.091059780849"
$ws.Range("C19").Value = [double]"0.6987839110824718"
$ws.Rows.Item(20).Delete()

# --- Sheet 7 ---
$ws = $wb.Worksheets.Item(7)
$ws.Name = "summ39426374"
$ws.Range("B2").Value = [double]"9770.310933659608"
$ws.Range("C2").Value = [double]"0.001759461015690735"
$ws.Range("B3").Value = [double]"256.3173113213695"
$ws.Range("C3").Value = [double]"0.8268493431636506"
$ws.Range("B4").Value = [double]"-1023.544345252612"
$ws.Range("C4").Value = [double]"0.4259102267716792"
$ws.Range("B5").Value = [double]"1008.125006794192"
$ws.Range("C5").Value = [double]"0.01883319584370745"
$ws.Range("B6").Value = [double]"-258.5751065024388"
$ws.Range("C6").Value = [double]"0.5574952883000601"
$ws.Range("B7").Value = [double]"794.8765827546322"
$ws.Range("C7").Value = [double]"0.1167479052300814"
$ws.Range("B8").Value = [double]"-146.5603205743207"
$ws.Range("C8").Value = [double]"0.7423337791203251"
$ws.Range("B9").Value = [double]"-146.9929382586647"
$ws.Range("C9").Value = [double]"0.3166693479930432"
$ws.Range("B10").Value = [double]"-1636.715659804823"
$ws.Range("C10").Value = [double]"1.000354866079533e-06"
$ws.Range("B11").Value = [double]"-64.52870804692766"
$ws.Range("C11").Value = [double]"3.903745922540236e-05"
$ws.Range("B12").Value = [double]"349.4114468838819"
$ws.Range("C12").Value = [double]"0.02007424702954179"
$ws.Range("B13").Value = [double]"520.9613393295032"
$ws.Range("C13").Value = [double]"4.235084025646359e-09"
$ws.Range("B14").Value = [double]"0.0235212941913644"
$ws.Range("C14").Value = [double]"0.7810304134159117"
$ws.Range("B15").Value = [double]"0.000207765692984948"
$ws.Range("C15").Value = [double]"0.03084735018420729"
$ws.Range("B16").Value = [double]"-39.34352357819594"
$ws.Range("C16").Value = [double]"0.05721854940877898"
$ws.Range("B17").Value = [double]"-2.958414885230631"
$ws.Range("C17").Value = [double]"0.8904019408150774"
$ws.Range("B18").Value = [double]"-156.2859534830804"
$ws.Range("C18").Value = [double]"0.92609283976093"
$ws.Range("B19").Value = [double]"1091.245169081173"
$ws.Range("C19").Value = [double]"0.5917723584914238"
$ws.Rows.Item(20).Delete()

# --- Sheet 8 ---
$ws = $wb.Worksheets.Item(8)
$ws.Name = "summ39773159"
$ws.Range("B2").Value = [double]"10703.34334497797"
$ws.Range("C2").Value = [double]"0.0007213345130088152"
$ws.Range("B3").Value = [double]"41.58699766453648"
$ws.Range("C3").Value = [double]"0.971834443522112"
$ws.Range("B4").Value = [double]"-353.2091545264382"
$ws.Range("C4").Value = [double]"0.7942445049913501"
$ws.Range("B5").Value = [double]"1060.499608487731"
$ws.Range("C5").Value = [double]"0.0161048635591307"
$ws.Range("B6").Value = [double]"66.20113780488637"
$ws.Range("C6").Value = [double]"0.8815631589298416"
$ws.Range("B7").Value = [double]"748.8854199734858"
$ws.Range("C7").Value = [double]"0.145921322007337"
$ws.Range("B8").Value = [double]"290.4609857936693"
$ws.Range("C8").Value = [double]"0.5165827053990747"
$ws.Range("B9").Value = [double]"-129.842673710341"
$ws.Range("C9").Value = [double]"0.377586657399763"
$ws.Range("B10").Value = [double]"-1662.772076826079"
$ws.Range("C10").Value = [double]"8.117904402065831e-07"
$ws.Range("B11").Value = [double]"-69.23968668050712"
$ws.Range("C11").Value = [double]"1.074389350070923e-05"
$ws.Range("B12").Value = [double]"426.7581772162032"
$ws.Range("C12").Value = [double]"0.005171775984612609"
$ws.Range("B13").Value = [double]"516.6450975298853"
$ws.Range("C13").Value = [double]"8.309048402014176e-09"
$ws.Range("B14").Value = [double]"-0.004317913348037258"
$ws.Range("C14").Value = [double]"0.9606626481062862"
$ws.Range("B15").Value = [double]"0.0001950415539772989"
$ws.Range("C15").Value = [double]"0.05574690352573724"
$ws.Range("B16").Value = [double]"-36.80062744287117"
$ws.Range("C16").Value = [double]"0.07979954175896874"
$ws.Range("B17").Value = [double]"-10.07664214144115"
$ws.Range("C17").Value = [double]"0.6421820262927396"
$ws.Range("B18").Value = [double]"-649.8896370503406"
$ws.Range("C18").Value = [double]"0.7028395392677917"
$ws.Range("B19").Value = [double]"651.1485249645776"
$ws.Range("C19").Value = [double]"0.7511361222649808"
$ws.Rows.Item(20).Delete()

# --- Sheet 9 ---
$ws = $wb.Worksheets.Item(9)
$ws.Name = "summ40026733"
$ws.Range("B2").Value = [double]"8632.320225368694"
$ws.Range("C2").Value = [double]"0.006832451176023911"
$ws.Range("B3").Value = [double]"252.8473429232644"
$ws.Range("C3").Value = [double]"0.8256534893664"
$ws.Range("B4").Value = [double]"-67.95150030913567"
$ws.Range("C4").Value = [double]"0.9577415456805283"
$ws.Range("B5").Value = [double]"1510.505877101646"
$ws.Range("C5").Value = [double]"0.0005430485613399295"
$ws.Range("B6").Value = [double]"-271.2526735550483"
$ws.Range("C6").Value = [double]"0.5420735129202412"
$ws.Range("B7").Value = [double]"555.7000390958182"
$ws.Range("C7").Value = [double]"0.2786204193473465"
$ws.Range("B8").Value = [double]"-70.10818449087631"
$ws.Range("C8").Value = [double]"0.8759178441451195"
$ws.Range("B9").Value = [double]"-228.4488193879441"
$ws.Range("C9").Value = [double]"0.123971319861603"
$ws.Range("B10").Value = [double]"-1745.404557972164"
$ws.Range("C10").Value = [double]"2.204512074290137e-07"
$ws.Range("B11").Value = [double]"-69.24042079037694"
$ws.Range("C11").Value = [double]"1.03548405952829e-05"
$ws.Range("B12").Value = [double]"413.1347591859451"
$ws.Range("C12").Value = [double]"0.006422036430735328"
$ws.Range("B13").Value = [double]"509.1145036074305"
$ws.Range("C13").Value = [double]"1.965232287988443e-08"
$ws.Range("B14").Value = [double]"0.0277740761972934"
$ws.Range("C14").Value = [double]"0.744560093119673"
$ws.Range("B15").Value = [double]"0.0001253327103065678"
$ws.Range("C15").Value = [double]"0.2093951782313385"
$ws.Range("B16").Value = [double]"-5.240150594966973"
$ws.Range("C16").Value = [double]"0.8032694331935271"
$ws.Range("B17").Value = [double]"6.836680893224445"
$ws.Range("C17").Value = [double]"0.7534920230983817"
$ws.Range("B18").Value = [double]"-738.6580968594246"
$ws.Range("C18").Value = [double]"0.6646691847383532"
$ws.Range("B19").Value = [double]"-1659.000307947287"
$ws.Range("C19").Value = [double]"0.4178417976619287"
$ws.Rows.Item(20).Delete()
